# ============================================================
# Kazakhstan Premier League - odds feed refresh (07-04-2024 22:30)
# For same-kickoff-time fixtures, the re-fetched feed returned the
# two matches in swapped order; realign row data (cols B:AC) while
# keeping the existing row index (col A) in place.
# ============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# swap row 9 and row 10 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(9,2).Value = 6221789
$ws.Cells.Item(10,2).Value = 6221694
$ws.Cells.Item(9,6).Value = 'Kairat Almaty'
$ws.Cells.Item(10,6).Value = 'FC Astana'
$ws.Cells.Item(9,7).Value = 'FK Kyzylzhar'
$ws.Cells.Item(10,7).Value = 'FK Kaspyi Aktau'
$ws.Cells.Item(9,9).Value = 0
$ws.Cells.Item(10,9).Value = 1
$ws.Cells.Item(9,11).Value = 2.25
$ws.Cells.Item(10,11).Value = 1.333
$ws.Cells.Item(9,12).Value = 3.2
$ws.Cells.Item(10,12).Value = 4.333
$ws.Cells.Item(9,13).Value = 2.8
$ws.Cells.Item(10,13).Value = 7.5
$ws.Cells.Item(9,14).Value = 2.1
$ws.Cells.Item(10,14).Value = 1.2
$ws.Cells.Item(9,15).Value = 3
$ws.Cells.Item(10,15).Value = 5.5
$ws.Cells.Item(9,16).Value = 3.2
$ws.Cells.Item(10,16).Value = 11
$ws.Cells.Item(9,17).Value = -0.25
$ws.Cells.Item(10,17).Value = -2
$ws.Cells.Item(9,18).Value = 1.875
$ws.Cells.Item(10,18).Value = 1.975
$ws.Cells.Item(9,19).Value = 1.925
$ws.Cells.Item(10,19).Value = 1.825
$ws.Cells.Item(9,20).Value = 2.25
$ws.Cells.Item(10,20).Value = 3
$ws.Cells.Item(9,21).Value = 2
$ws.Cells.Item(10,21).Value = 1.95
$ws.Cells.Item(9,22).Value = 1.8
$ws.Cells.Item(10,22).Value = 1.85
$ws.Cells.Item(9,23).Value = 1.1
$ws.Cells.Item(10,23).Value = 0.2
$ws.Cells.Item(9,26).Value = 0.875
$ws.Cells.Item(10,26).Value = 0
$ws.Cells.Item(9,27).Value = -1
$ws.Cells.Item(10,27).Value = -0
$ws.Cells.Item(9,28).Value = 1
$ws.Cells.Item(10,28).Value = 0.95

# swap row 16 and row 17 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(16,2).Value = 6221698
$ws.Cells.Item(17,2).Value = 6221693
$ws.Cells.Item(16,6).Value = 'FK Maktaaral'
$ws.Cells.Item(17,6).Value = 'Zhetysu'
$ws.Cells.Item(16,7).Value = 'FK Aktobe'
$ws.Cells.Item(17,7).Value = 'Shakhter Karagandy'
$ws.Cells.Item(16,9).Value = 2
$ws.Cells.Item(17,9).Value = 3
$ws.Cells.Item(16,11).Value = 4.333
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(16,12).Value = 3.5
$ws.Cells.Item(17,12).Value = 3.4
$ws.Cells.Item(16,13).Value = 1.666
$ws.Cells.Item(17,13).Value = 3.1
$ws.Cells.Item(16,14).Value = 4.2
$ws.Cells.Item(17,14).Value = 2.2
$ws.Cells.Item(16,15).Value = 3.4
$ws.Cells.Item(17,15).Value = 3.3
$ws.Cells.Item(16,16).Value = 1.7
$ws.Cells.Item(17,16).Value = 2.8
$ws.Cells.Item(16,17).Value = 0.75
$ws.Cells.Item(17,17).Value = -0.25
$ws.Cells.Item(16,18).Value = 1.825
$ws.Cells.Item(17,18).Value = 1.95
$ws.Cells.Item(16,19).Value = 1.975
$ws.Cells.Item(17,19).Value = 1.85
$ws.Cells.Item(16,21).Value = 1.925
$ws.Cells.Item(17,21).Value = 1.85
$ws.Cells.Item(16,22).Value = 1.875
$ws.Cells.Item(17,22).Value = 1.95
$ws.Cells.Item(16,25).Value = 0.7
$ws.Cells.Item(17,25).Value = 1.8
$ws.Cells.Item(16,26).Value = -0.5
$ws.Cells.Item(17,26).Value = -1
$ws.Cells.Item(16,27).Value = 0.4875
$ws.Cells.Item(17,27).Value = 0.8500000000000001
$ws.Cells.Item(16,28).Value = 0.925
$ws.Cells.Item(17,28).Value = 0.8500000000000001

# swap row 37 and row 38 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(37,2).Value = 6221712
$ws.Cells.Item(38,2).Value = 6221708
$ws.Cells.Item(37,6).Value = 'FK Aksu'
$ws.Cells.Item(38,6).Value = 'Kaisar Kyzylorda'
$ws.Cells.Item(37,7).Value = 'Shakhter Karagandy'
$ws.Cells.Item(38,7).Value = 'Kairat Almaty'
$ws.Cells.Item(37,8).Value = 2
$ws.Cells.Item(38,8).Value = 0
$ws.Cells.Item(37,9).Value = 1
$ws.Cells.Item(38,9).Value = 0
$ws.Cells.Item(37,10).Value = 'H'
$ws.Cells.Item(38,10).Value = 'D'
$ws.Cells.Item(37,11).Value = 2.1
$ws.Cells.Item(38,11).Value = 3
$ws.Cells.Item(37,12).Value = 3.25
$ws.Cells.Item(38,12).Value = 3.4
$ws.Cells.Item(37,13).Value = 3
$ws.Cells.Item(38,13).Value = 2.05
$ws.Cells.Item(37,14).Value = 2.15
$ws.Cells.Item(38,14).Value = 3.2
$ws.Cells.Item(37,15).Value = 3.25
$ws.Cells.Item(38,15).Value = 3.4
$ws.Cells.Item(37,16).Value = 2.9
$ws.Cells.Item(38,16).Value = 1.95
$ws.Cells.Item(37,17).Value = -0.25
$ws.Cells.Item(38,17).Value = 0.5
$ws.Cells.Item(37,18).Value = 1.95
$ws.Cells.Item(38,18).Value = 1.75
$ws.Cells.Item(37,19).Value = 1.85
$ws.Cells.Item(38,19).Value = 1.95
$ws.Cells.Item(37,20).Value = 2.5
$ws.Cells.Item(38,20).Value = 2.25
$ws.Cells.Item(37,21).Value = 1.975
$ws.Cells.Item(38,21).Value = 1.925
$ws.Cells.Item(37,22).Value = 1.825
$ws.Cells.Item(38,22).Value = 1.875
$ws.Cells.Item(37,23).Value = 1.15
$ws.Cells.Item(38,23).Value = -1
$ws.Cells.Item(37,24).Value = -1
$ws.Cells.Item(38,24).Value = 2.4
$ws.Cells.Item(37,26).Value = 0.95
$ws.Cells.Item(38,26).Value = 0.75
$ws.Cells.Item(37,28).Value = 0.9750000000000001
$ws.Cells.Item(38,28).Value = -1
$ws.Cells.Item(37,29).Value = -1
$ws.Cells.Item(38,29).Value = 0.875

# swap row 50 and row 51 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(50,2).Value = 7055064
$ws.Cells.Item(51,2).Value = 6221723
$ws.Cells.Item(50,6).Value = 'Kaisar Kyzylorda'
$ws.Cells.Item(51,6).Value = 'FK Aksu'
$ws.Cells.Item(50,7).Value = 'Ordabasy'
$ws.Cells.Item(51,7).Value = 'FK Aktobe'
$ws.Cells.Item(50,8).Value = 2
$ws.Cells.Item(51,8).Value = 1
$ws.Cells.Item(50,9).Value = 3
$ws.Cells.Item(51,9).Value = 1
$ws.Cells.Item(50,10).Value = 'A'
$ws.Cells.Item(51,10).Value = 'D'
$ws.Cells.Item(50,11).Value = 6.5
$ws.Cells.Item(51,11).Value = 3.75
$ws.Cells.Item(50,12).Value = 4.5
$ws.Cells.Item(51,12).Value = 3.3
$ws.Cells.Item(50,13).Value = 1.363
$ws.Cells.Item(51,13).Value = 1.833
$ws.Cells.Item(50,14).Value = 4.2
$ws.Cells.Item(51,14).Value = 2.6
$ws.Cells.Item(50,15).Value = 4
$ws.Cells.Item(51,15).Value = 3.1
$ws.Cells.Item(50,16).Value = 1.6
$ws.Cells.Item(51,16).Value = 2.5
$ws.Cells.Item(50,17).Value = 1
$ws.Cells.Item(51,17).Value = 0
$ws.Cells.Item(50,18).Value = 1.725
$ws.Cells.Item(51,18).Value = 1.925
$ws.Cells.Item(50,19).Value = 2.075
$ws.Cells.Item(51,19).Value = 1.875
$ws.Cells.Item(50,20).Value = 2.25
$ws.Cells.Item(51,20).Value = 2.5
$ws.Cells.Item(50,21).Value = 1.875
$ws.Cells.Item(51,21).Value = 1.9
$ws.Cells.Item(50,22).Value = 1.925
$ws.Cells.Item(51,22).Value = 1.9
$ws.Cells.Item(50,24).Value = -1
$ws.Cells.Item(51,24).Value = 2.1
$ws.Cells.Item(50,25).Value = 0.6000000000000001
$ws.Cells.Item(51,25).Value = -1
$ws.Cells.Item(50,28).Value = 0.875
$ws.Cells.Item(51,28).Value = -1
$ws.Cells.Item(50,29).Value = -1
$ws.Cells.Item(51,29).Value = 0.8999999999999999

# swap row 63 and row 64 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(63,2).Value = 6221729
$ws.Cells.Item(64,2).Value = 6221732
$ws.Cells.Item(63,6).Value = 'Kairat Almaty'
$ws.Cells.Item(64,6).Value = 'FK Atyrau'
$ws.Cells.Item(63,7).Value = 'FC Astana'
$ws.Cells.Item(64,7).Value = 'FK Aktobe'
$ws.Cells.Item(63,8).Value = 1
$ws.Cells.Item(64,8).Value = 0
$ws.Cells.Item(63,10).Value = 'H'
$ws.Cells.Item(64,10).Value = 'D'
$ws.Cells.Item(63,11).Value = 2.9
$ws.Cells.Item(64,11).Value = 3.3
$ws.Cells.Item(63,12).Value = 3.3
$ws.Cells.Item(64,12).Value = 3.25
$ws.Cells.Item(63,13).Value = 2.15
$ws.Cells.Item(64,13).Value = 2
$ws.Cells.Item(63,14).Value = 1.75
$ws.Cells.Item(64,14).Value = 3.5
$ws.Cells.Item(63,15).Value = 3.5
$ws.Cells.Item(64,15).Value = 3.25
$ws.Cells.Item(63,16).Value = 4
$ws.Cells.Item(64,16).Value = 1.909
$ws.Cells.Item(63,17).Value = -0.5
$ws.Cells.Item(64,17).Value = 0.5
$ws.Cells.Item(63,18).Value = 1.8
$ws.Cells.Item(64,18).Value = 1.775
$ws.Cells.Item(63,19).Value = 2
$ws.Cells.Item(64,19).Value = 2.025
$ws.Cells.Item(63,20).Value = 2.5
$ws.Cells.Item(64,20).Value = 2.25
$ws.Cells.Item(63,23).Value = 0.75
$ws.Cells.Item(64,23).Value = -1
$ws.Cells.Item(63,24).Value = -1
$ws.Cells.Item(64,24).Value = 2.25
$ws.Cells.Item(63,26).Value = 0.8
$ws.Cells.Item(64,26).Value = 0.7749999999999999

# swap row 85 and row 86 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(85,2).Value = 6221743
$ws.Cells.Item(86,2).Value = 6221809
$ws.Cells.Item(85,6).Value = 'FK Aksu'
$ws.Cells.Item(86,6).Value = 'FK Kyzylzhar'
$ws.Cells.Item(85,7).Value = 'Kaisar Kyzylorda'
$ws.Cells.Item(86,7).Value = 'Okzhetpes Kokshetau'
$ws.Cells.Item(85,8).Value = 2
$ws.Cells.Item(86,8).Value = 1
$ws.Cells.Item(85,9).Value = 0
$ws.Cells.Item(86,9).Value = 2
$ws.Cells.Item(85,10).Value = 'H'
$ws.Cells.Item(86,10).Value = 'A'
$ws.Cells.Item(85,11).Value = 3
$ws.Cells.Item(86,11).Value = 2.05
$ws.Cells.Item(85,13).Value = 2.2
$ws.Cells.Item(86,13).Value = 3.3
$ws.Cells.Item(85,14).Value = 2.9
$ws.Cells.Item(86,14).Value = 1.55
$ws.Cells.Item(85,15).Value = 3.1
$ws.Cells.Item(86,15).Value = 3.5
$ws.Cells.Item(85,16).Value = 2.25
$ws.Cells.Item(86,16).Value = 5.25
$ws.Cells.Item(85,17).Value = 0.25
$ws.Cells.Item(86,17).Value = -1
$ws.Cells.Item(85,18).Value = 1.8
$ws.Cells.Item(86,18).Value = 2
$ws.Cells.Item(85,19).Value = 2
$ws.Cells.Item(86,19).Value = 1.8
$ws.Cells.Item(85,21).Value = 1.95
$ws.Cells.Item(86,21).Value = 1.875
$ws.Cells.Item(85,22).Value = 1.85
$ws.Cells.Item(86,22).Value = 1.925
$ws.Cells.Item(85,23).Value = 1.9
$ws.Cells.Item(86,23).Value = -1
$ws.Cells.Item(85,25).Value = -1
$ws.Cells.Item(86,25).Value = 4.25
$ws.Cells.Item(85,26).Value = 0.8
$ws.Cells.Item(86,26).Value = -1
$ws.Cells.Item(85,27).Value = -1
$ws.Cells.Item(86,27).Value = 0.8
$ws.Cells.Item(85,28).Value = -0.5
$ws.Cells.Item(86,28).Value = 0.875
$ws.Cells.Item(85,29).Value = 0.425
$ws.Cells.Item(86,29).Value = -1

# swap row 92 and row 93 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(92,2).Value = 6221748
$ws.Cells.Item(93,2).Value = 6221749
$ws.Cells.Item(92,6).Value = 'Tobol Kostanay'
$ws.Cells.Item(93,6).Value = 'Kairat Almaty'
$ws.Cells.Item(92,7).Value = 'FK Kaspyi Aktau'
$ws.Cells.Item(93,7).Value = 'FK Aksu'
$ws.Cells.Item(92,8).Value = 1
$ws.Cells.Item(93,8).Value = 4
$ws.Cells.Item(92,10).Value = 'D'
$ws.Cells.Item(93,10).Value = 'H'
$ws.Cells.Item(92,11).Value = 1.533
$ws.Cells.Item(93,11).Value = 1.5
$ws.Cells.Item(92,13).Value = 6
$ws.Cells.Item(93,13).Value = 6.5
$ws.Cells.Item(92,14).Value = 1.444
$ws.Cells.Item(93,14).Value = 1.333
$ws.Cells.Item(92,15).Value = 4.2
$ws.Cells.Item(93,15).Value = 4.5
$ws.Cells.Item(92,16).Value = 6.5
$ws.Cells.Item(93,16).Value = 7.5
$ws.Cells.Item(92,17).Value = -1.25
$ws.Cells.Item(93,17).Value = -1.5
$ws.Cells.Item(92,18).Value = 1.85
$ws.Cells.Item(93,18).Value = 1.8
$ws.Cells.Item(92,19).Value = 1.95
$ws.Cells.Item(93,19).Value = 2
$ws.Cells.Item(92,21).Value = 1.925
$ws.Cells.Item(93,21).Value = 1.85
$ws.Cells.Item(92,22).Value = 1.875
$ws.Cells.Item(93,22).Value = 1.95
$ws.Cells.Item(92,23).Value = -1
$ws.Cells.Item(93,23).Value = 0.333
$ws.Cells.Item(92,24).Value = 3.2
$ws.Cells.Item(93,24).Value = -1
$ws.Cells.Item(92,26).Value = -1
$ws.Cells.Item(93,26).Value = 0.8
$ws.Cells.Item(92,27).Value = 0.95
$ws.Cells.Item(93,27).Value = -1
$ws.Cells.Item(92,28).Value = -1
$ws.Cells.Item(93,28).Value = 0.8500000000000001
$ws.Cells.Item(92,29).Value = 0.875
$ws.Cells.Item(93,29).Value = -1

# swap row 99 and row 100 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(99,2).Value = 6221815
$ws.Cells.Item(100,2).Value = 6221752
$ws.Cells.Item(99,6).Value = 'FK Atyrau'
$ws.Cells.Item(100,6).Value = 'FK Kyzylzhar'
$ws.Cells.Item(99,7).Value = 'Kairat Almaty'
$ws.Cells.Item(100,7).Value = 'Kaisar Kyzylorda'
$ws.Cells.Item(99,9).Value = 0
$ws.Cells.Item(100,9).Value = 1
$ws.Cells.Item(99,10).Value = 'D'
$ws.Cells.Item(100,10).Value = 'A'
$ws.Cells.Item(99,11).Value = 3
$ws.Cells.Item(100,11).Value = 1.833
$ws.Cells.Item(99,12).Value = 3
$ws.Cells.Item(100,12).Value = 3.2
$ws.Cells.Item(99,13).Value = 2.25
$ws.Cells.Item(100,13).Value = 4
$ws.Cells.Item(99,14).Value = 3.1
$ws.Cells.Item(100,14).Value = 1.85
$ws.Cells.Item(99,15).Value = 3.1
$ws.Cells.Item(100,15).Value = 3.2
$ws.Cells.Item(99,16).Value = 2.15
$ws.Cells.Item(100,16).Value = 4
$ws.Cells.Item(99,17).Value = 0.25
$ws.Cells.Item(100,17).Value = -0.5
$ws.Cells.Item(99,18).Value = 1.85
$ws.Cells.Item(100,18).Value = 1.9
$ws.Cells.Item(99,19).Value = 1.95
$ws.Cells.Item(100,19).Value = 1.9
$ws.Cells.Item(99,20).Value = 2.25
$ws.Cells.Item(100,20).Value = 2
$ws.Cells.Item(99,21).Value = 1.8
$ws.Cells.Item(100,21).Value = 1.775
$ws.Cells.Item(99,22).Value = 2
$ws.Cells.Item(100,22).Value = 2.025
$ws.Cells.Item(99,24).Value = 2.1
$ws.Cells.Item(100,24).Value = -1
$ws.Cells.Item(99,25).Value = -1
$ws.Cells.Item(100,25).Value = 3
$ws.Cells.Item(99,26).Value = 0.425
$ws.Cells.Item(100,26).Value = -1
$ws.Cells.Item(99,27).Value = -0.5
$ws.Cells.Item(100,27).Value = 0.8999999999999999
$ws.Cells.Item(99,29).Value = 1
$ws.Cells.Item(100,29).Value = 1.025

# swap row 101 and row 104 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(101,2).Value = 6221816
$ws.Cells.Item(104,2).Value = 6221814
$ws.Cells.Item(101,6).Value = 'FK Kaspyi Aktau'
$ws.Cells.Item(104,6).Value = 'Okzhetpes Kokshetau'
$ws.Cells.Item(101,7).Value = 'Ordabasy'
$ws.Cells.Item(104,7).Value = 'FK Maktaaral'
$ws.Cells.Item(101,9).Value = 0
$ws.Cells.Item(104,9).Value = 1
$ws.Cells.Item(101,10).Value = 'H'
$ws.Cells.Item(104,10).Value = 'D'
$ws.Cells.Item(101,11).Value = 3.4
$ws.Cells.Item(104,11).Value = 2.3
$ws.Cells.Item(101,12).Value = 3.4
$ws.Cells.Item(104,12).Value = 3.1
$ws.Cells.Item(101,13).Value = 1.909
$ws.Cells.Item(104,13).Value = 2.8
$ws.Cells.Item(101,14).Value = 4.2
$ws.Cells.Item(104,14).Value = 2.3
$ws.Cells.Item(101,15).Value = 4
$ws.Cells.Item(104,15).Value = 3.1
$ws.Cells.Item(101,16).Value = 1.571
$ws.Cells.Item(104,16).Value = 2.8
$ws.Cells.Item(101,17).Value = 0.75
$ws.Cells.Item(104,17).Value = 0
$ws.Cells.Item(101,18).Value = 1.95
$ws.Cells.Item(104,18).Value = 1.75
$ws.Cells.Item(101,19).Value = 1.75
$ws.Cells.Item(104,19).Value = 2.05
$ws.Cells.Item(101,20).Value = 3
$ws.Cells.Item(104,20).Value = 2.25
$ws.Cells.Item(101,21).Value = 1.975
$ws.Cells.Item(104,21).Value = 1.875
$ws.Cells.Item(101,22).Value = 1.825
$ws.Cells.Item(104,22).Value = 1.925
$ws.Cells.Item(101,23).Value = 3.2
$ws.Cells.Item(104,23).Value = -1
$ws.Cells.Item(101,24).Value = -1
$ws.Cells.Item(104,24).Value = 2.1
$ws.Cells.Item(101,26).Value = 0.95
$ws.Cells.Item(104,26).Value = 0
$ws.Cells.Item(101,27).Value = -1
$ws.Cells.Item(104,27).Value = -0
$ws.Cells.Item(101,28).Value = -1
$ws.Cells.Item(104,28).Value = -0.5
$ws.Cells.Item(101,29).Value = 0.825
$ws.Cells.Item(104,29).Value = 0.4625

# swap row 119 and row 120 (re-ordered fixtures at identical kickoff time)
$ws.Cells.Item(119,2).Value = 7873759
$ws.Cells.Item(120,2).Value = 7874795
$ws.Cells.Item(119,6).Value = 'FK Aktobe'
$ws.Cells.Item(120,6).Value = 'FK Kyzylzhar'
$ws.Cells.Item(119,7).Value = 'FK Zhenys'
$ws.Cells.Item(120,7).Value = 'Tobol Kostanay'
$ws.Cells.Item(119,8).Value = 3
$ws.Cells.Item(120,8).Value = 0
$ws.Cells.Item(119,10).Value = 'H'
$ws.Cells.Item(120,10).Value = 'D'
$ws.Cells.Item(119,11).Value = 1.25
$ws.Cells.Item(120,11).Value = 2.2
$ws.Cells.Item(119,12).Value = 5.75
$ws.Cells.Item(120,12).Value = 3.1
$ws.Cells.Item(119,13).Value = 7
$ws.Cells.Item(120,13).Value = 3
$ws.Cells.Item(119,14).Value = 1.444
$ws.Cells.Item(120,14).Value = 2.625
$ws.Cells.Item(119,15).Value = 4.75
$ws.Cells.Item(120,15).Value = 3
$ws.Cells.Item(119,16).Value = 4.75
$ws.Cells.Item(120,16).Value = 2.55
$ws.Cells.Item(119,17).Value = -1.25
$ws.Cells.Item(120,17).Value = 0
$ws.Cells.Item(119,18).Value = 1.95
$ws.Cells.Item(120,18).Value = 1.9
$ws.Cells.Item(119,19).Value = 1.85
$ws.Cells.Item(120,19).Value = 1.9
$ws.Cells.Item(119,20).Value = 2.75
$ws.Cells.Item(120,20).Value = 2
$ws.Cells.Item(119,21).Value = 1.9
$ws.Cells.Item(120,21).Value = 1.95
$ws.Cells.Item(119,22).Value = 1.9
$ws.Cells.Item(120,22).Value = 1.85
$ws.Cells.Item(119,23).Value = 0.444
$ws.Cells.Item(120,23).Value = -1
$ws.Cells.Item(119,24).Value = -1
$ws.Cells.Item(120,24).Value = 2
$ws.Cells.Item(119,26).Value = 0.95
$ws.Cells.Item(120,26).Value = 0
$ws.Cells.Item(119,27).Value = -1
$ws.Cells.Item(120,27).Value = -0
$ws.Cells.Item(119,28).Value = 0.45
$ws.Cells.Item(120,28).Value = -1
$ws.Cells.Item(119,29).Value = -0.5
$ws.Cells.Item(120,29).Value = 0.8500000000000001

# ------------------------------------------------------------
# Append 6 new fixtures (rows 123-128) to the bottom of the table
# ------------------------------------------------------------
# Column A (bold/bordered/centered) and column E (date format)
# reuse the style already applied on row 122 via copy/paste-format
# so no duplicate style entries are added to styles.xml.
$ws.Cells.Item(122,1).Copy() | Out-Null
$ws.Cells.Item(123,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,5).Copy() | Out-Null
$ws.Cells.Item(123,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,1).Copy() | Out-Null
$ws.Cells.Item(124,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,5).Copy() | Out-Null
$ws.Cells.Item(124,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,1).Copy() | Out-Null
$ws.Cells.Item(125,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,5).Copy() | Out-Null
$ws.Cells.Item(125,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,1).Copy() | Out-Null
$ws.Cells.Item(126,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,5).Copy() | Out-Null
$ws.Cells.Item(126,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,1).Copy() | Out-Null
$ws.Cells.Item(127,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,5).Copy() | Out-Null
$ws.Cells.Item(127,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,1).Copy() | Out-Null
$ws.Cells.Item(128,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122,5).Copy() | Out-Null
$ws.Cells.Item(128,5).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# row 123
$ws.Cells.Item(123,3).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(123,4).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(123,1).Value = 121
$ws.Cells.Item(123,2).Value = 7873757
$ws.Cells.Item(123,5).Value = 45388.25
$ws.Cells.Item(123,6).Value = 'FK Kyzylzhar'
$ws.Cells.Item(123,7).Value = 'FK Zhenys'
$ws.Cells.Item(123,8).Value = 4
$ws.Cells.Item(123,9).Value = 0
$ws.Cells.Item(123,10).Value = 'H'
$ws.Cells.Item(123,11).Value = 2
$ws.Cells.Item(123,12).Value = 3.25
$ws.Cells.Item(123,13).Value = 3.25
$ws.Cells.Item(123,14).Value = 1.7
$ws.Cells.Item(123,15).Value = 3.3
$ws.Cells.Item(123,16).Value = 4.75
$ws.Cells.Item(123,17).Value = -0.75
$ws.Cells.Item(123,18).Value = 1.95
$ws.Cells.Item(123,19).Value = 1.85
$ws.Cells.Item(123,20).Value = 2.25
$ws.Cells.Item(123,21).Value = 2
$ws.Cells.Item(123,22).Value = 1.8
$ws.Cells.Item(123,23).Value = 0.7
$ws.Cells.Item(123,24).Value = -1
$ws.Cells.Item(123,25).Value = -1
$ws.Cells.Item(123,26).Value = 0.95
$ws.Cells.Item(123,27).Value = -1
$ws.Cells.Item(123,28).Value = 1
$ws.Cells.Item(123,29).Value = -1

# row 124
$ws.Cells.Item(124,3).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(124,4).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(124,1).Value = 122
$ws.Cells.Item(124,2).Value = 7874797
$ws.Cells.Item(124,5).Value = 45388.35416666666
$ws.Cells.Item(124,6).Value = 'FC Elimai Semey'
$ws.Cells.Item(124,7).Value = 'Kaisar Kyzylorda'
$ws.Cells.Item(124,8).Value = 1
$ws.Cells.Item(124,9).Value = 0
$ws.Cells.Item(124,10).Value = 'H'
$ws.Cells.Item(124,11).Value = 2.2
$ws.Cells.Item(124,12).Value = 3.5
$ws.Cells.Item(124,13).Value = 2.7
$ws.Cells.Item(124,14).Value = 1.6
$ws.Cells.Item(124,15).Value = 3.8
$ws.Cells.Item(124,16).Value = 4.333
$ws.Cells.Item(124,17).Value = -0.75
$ws.Cells.Item(124,18).Value = 1.8
$ws.Cells.Item(124,19).Value = 2
$ws.Cells.Item(124,20).Value = 2.25
$ws.Cells.Item(124,21).Value = 1.85
$ws.Cells.Item(124,22).Value = 1.95
$ws.Cells.Item(124,23).Value = 0.6000000000000001
$ws.Cells.Item(124,24).Value = -1
$ws.Cells.Item(124,25).Value = -1
$ws.Cells.Item(124,26).Value = 0.4
$ws.Cells.Item(124,27).Value = -0.5
$ws.Cells.Item(124,28).Value = -1
$ws.Cells.Item(124,29).Value = 0.95

# row 125
$ws.Cells.Item(125,3).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(125,4).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(125,1).Value = 123
$ws.Cells.Item(125,2).Value = 7874798
$ws.Cells.Item(125,5).Value = 45388.45833333334
$ws.Cells.Item(125,6).Value = 'Kairat Almaty'
$ws.Cells.Item(125,7).Value = 'Zhetysu'
$ws.Cells.Item(125,8).Value = 0
$ws.Cells.Item(125,9).Value = 0
$ws.Cells.Item(125,10).Value = 'D'
$ws.Cells.Item(125,11).Value = 1.727
$ws.Cells.Item(125,12).Value = 3.75
$ws.Cells.Item(125,13).Value = 3.75
$ws.Cells.Item(125,14).Value = 1.571
$ws.Cells.Item(125,15).Value = 3.75
$ws.Cells.Item(125,16).Value = 4.75
$ws.Cells.Item(125,17).Value = -1
$ws.Cells.Item(125,18).Value = 1.8
$ws.Cells.Item(125,19).Value = 2
$ws.Cells.Item(125,20).Value = 2.5
$ws.Cells.Item(125,21).Value = 1.95
$ws.Cells.Item(125,22).Value = 1.75
$ws.Cells.Item(125,23).Value = -1
$ws.Cells.Item(125,24).Value = 2.75
$ws.Cells.Item(125,25).Value = -1
$ws.Cells.Item(125,26).Value = -1
$ws.Cells.Item(125,27).Value = 1
$ws.Cells.Item(125,28).Value = -1
$ws.Cells.Item(125,29).Value = 0.75

# row 126
$ws.Cells.Item(126,3).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(126,4).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(126,1).Value = 124
$ws.Cells.Item(126,2).Value = 7874799
$ws.Cells.Item(126,5).Value = 45389.25
$ws.Cells.Item(126,6).Value = 'FK Aktobe'
$ws.Cells.Item(126,7).Value = 'FK Atyrau'
$ws.Cells.Item(126,8).Value = 1
$ws.Cells.Item(126,9).Value = 1
$ws.Cells.Item(126,10).Value = 'D'
$ws.Cells.Item(126,11).Value = 1.85
$ws.Cells.Item(126,12).Value = 3.25
$ws.Cells.Item(126,13).Value = 3.75
$ws.Cells.Item(126,14).Value = 1.65
$ws.Cells.Item(126,15).Value = 3.4
$ws.Cells.Item(126,16).Value = 4.5
$ws.Cells.Item(126,17).Value = -0.75
$ws.Cells.Item(126,18).Value = 1.9
$ws.Cells.Item(126,19).Value = 1.9
$ws.Cells.Item(126,20).Value = 2.25
$ws.Cells.Item(126,21).Value = 2.025
$ws.Cells.Item(126,22).Value = 1.775
$ws.Cells.Item(126,23).Value = -1
$ws.Cells.Item(126,24).Value = 2.4
$ws.Cells.Item(126,25).Value = -1
$ws.Cells.Item(126,26).Value = -1
$ws.Cells.Item(126,27).Value = 0.8999999999999999
$ws.Cells.Item(126,28).Value = -0.5
$ws.Cells.Item(126,29).Value = 0.3875

# row 127
$ws.Cells.Item(127,3).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(127,4).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(127,1).Value = 125
$ws.Cells.Item(127,2).Value = 7871216
$ws.Cells.Item(127,5).Value = 45389.35416666666
$ws.Cells.Item(127,6).Value = 'Shakhter Karagandy'
$ws.Cells.Item(127,7).Value = 'FC Astana'
$ws.Cells.Item(127,8).Value = 1
$ws.Cells.Item(127,9).Value = 0
$ws.Cells.Item(127,10).Value = 'H'
$ws.Cells.Item(127,11).Value = 2.6
$ws.Cells.Item(127,12).Value = 3
$ws.Cells.Item(127,13).Value = 2.5
$ws.Cells.Item(127,14).Value = 6
$ws.Cells.Item(127,15).Value = 4
$ws.Cells.Item(127,16).Value = 1.4
$ws.Cells.Item(127,17).Value = 1.25
$ws.Cells.Item(127,18).Value = 1.8
$ws.Cells.Item(127,19).Value = 2
$ws.Cells.Item(127,20).Value = 2.25
$ws.Cells.Item(127,21).Value = 1.85
$ws.Cells.Item(127,22).Value = 1.95
$ws.Cells.Item(127,23).Value = 5
$ws.Cells.Item(127,24).Value = -1
$ws.Cells.Item(127,25).Value = -1
$ws.Cells.Item(127,26).Value = 0.8
$ws.Cells.Item(127,27).Value = -1
$ws.Cells.Item(127,28).Value = -1
$ws.Cells.Item(127,29).Value = 0.95

# row 128
$ws.Cells.Item(128,3).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(128,4).Value = 'Kazakhstan Premier League'
$ws.Cells.Item(128,1).Value = 126
$ws.Cells.Item(128,2).Value = 7874801
$ws.Cells.Item(128,5).Value = 45389.45833333334
$ws.Cells.Item(128,6).Value = 'Turan'
$ws.Cells.Item(128,7).Value = 'Ordabasy'
$ws.Cells.Item(128,8).Value = 0
$ws.Cells.Item(128,9).Value = 1
$ws.Cells.Item(128,10).Value = 'A'
$ws.Cells.Item(128,11).Value = 2.5
$ws.Cells.Item(128,12).Value = 3.75
$ws.Cells.Item(128,13).Value = 2.25
$ws.Cells.Item(128,14).Value = 5.75
$ws.Cells.Item(128,15).Value = 4.333
$ws.Cells.Item(128,16).Value = 1.444
$ws.Cells.Item(128,17).Value = 1.25
$ws.Cells.Item(128,18).Value = 1.775
$ws.Cells.Item(128,19).Value = 2.025
$ws.Cells.Item(128,20).Value = 2.25
$ws.Cells.Item(128,21).Value = 1.9
$ws.Cells.Item(128,22).Value = 1.9
$ws.Cells.Item(128,23).Value = -1
$ws.Cells.Item(128,24).Value = -1
$ws.Cells.Item(128,25).Value = 0.444
$ws.Cells.Item(128,26).Value = 0.3875
$ws.Cells.Item(128,27).Value = -0.5
$ws.Cells.Item(128,28).Value = -1
$ws.Cells.Item(128,29).Value = 0.8999999999999999
